# Added results for fine grained classification Config 6 Label Powerset Ubuntu
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ubuntu")

# Fill in the Config 6 (ngram(3) + POS) result rows for each model block
$ws.Range("C45").Value = "0.672 0.497 0.265 0.399 0.753"
$ws.Range("D45").Value = "0.578 0.548 0.282 0.469 0.763"
$ws.Range("E45").Value = "0.547 0.337 0.153 0.250 0.607"
$ws.Range("F45").Value = "0.790 0.861 0.916 0.944 0.976"
$ws.Range("C46").Value = "0.736 0.601 0.208 0.356 0.677"
$ws.Range("D46").Value = "0.470 0.484 0.333 0.440 0.822 "
$ws.Range("E46").Value = "0.750 0.446 0.116 0.217 0.513 "
$ws.Range("F46").Value = "0.730 0.848 0.924 0.942 0.975"
$ws.Range("C47").Value = "0.648 0.622 0.340 0.499 0.790"
$ws.Range("D47").Value = "0.584 0.492 0.271 0.447 0.700 "
$ws.Range("E47").Value = "0.512 0.471 0.206 0.336 0.658"
$ws.Range("F47").Value = "0.789 0.850 0.908 0.941 0.974"
$ws.Range("C48").Value = "0.492 0.205 0.062 0.039 0.713"
$ws.Range("D48").Value = "0.438 0.485 0.207 0.143 0.867"
$ws.Range("E48").Value = "0.345 0.115 0.032 0.020 0.556"
$ws.Range("F48").Value = "0.729 0.852 0.926 0.940 0.978"
$ws.Range("C49").Value = "0.639 0.356 0.191 0.282 0.755 "
$ws.Range("D49").Value = "0.619 0.881 0.952 0.781 0.947 "
$ws.Range("E49").Value = "0.495 0.217 0.106 0.164 0.607 "
$ws.Range("F49").Value = "0.800 0.880 0.939 0.952 0.982"
$ws.Range("C51").Value = "0.671 0.522 0.369 0.407 0.759 "
$ws.Range("D51").Value = "0.584 0.519 0.368 0.433 0.735"
$ws.Range("E51").Value = "0.544 0.361 0.228 0.257 0.615"
$ws.Range("F51").Value = "0.792 0.857 0.921 0.941 0.974"
$ws.Range("C52").Value = "0.752 0.588 0.217 0.364 0.713"
$ws.Range("D52").Value = "0.501 0.494 0.354 0.442 0.867"
$ws.Range("E52").Value = "0.747 0.432 0.122 0.224 0.556 "
$ws.Range("F52").Value = "0.754 0.851 0.925 0.942 0.978"
$ws.Range("C53").Value = "0.651 0.629 0.368 0.499 0.790 "
$ws.Range("D53").Value = "0.589 0.510 0.289 0.436 0.700"
$ws.Range("E53").Value = "0.515 0.478 0.228 0.336 0.658 "
$ws.Range("F53").Value = "0.792 0.855 0.909 0.940 0.974"
$ws.Range("C54").Value = "0.487 0.209 0.062 0.039 0.706"
$ws.Range("D54").Value = "0.432 0.471 0.207 0.143 0.865"
$ws.Range("E54").Value = "0.341 0.117 0.032 0.020 0.547"
$ws.Range("F54").Value = "0.726 0.850 0.926 0.940 0.977"
$ws.Range("C55").Value = "0.644 0.353 0.191 0.282 0.761"
$ws.Range("D55").Value = "0.635 0.854 0.952 0.758 0.947"
$ws.Range("E55").Value = "0.499 0.215 0.106 0.164 0.615 "
$ws.Range("F55").Value = "0.806 0.879 0.939 0.951 0.982"
$ws.Range("C57").Value = "0.690 0.538 0.348 0.373 0.772 "
$ws.Range("D57").Value = "0.602 0.546 0.377 0.417 0.771"
$ws.Range("E57").Value = "0.569 0.376 0.212 0.230 0.632 "
$ws.Range("F57").Value = "0.801 0.862 0.923 0.940 0.977"
$ws.Range("C58").Value = "0.761 0.533 0.165 0.169 0.707"
$ws.Range("D58").Value = "0.534 0.549 0.447 0.400 0.928"
$ws.Range("E58").Value = "0.734 0.371 0.090 0.092 0.547 "
$ws.Range("F58").Value = "0.776 0.862 0.931 0.943 0.979"
$ws.Range("C59").Value = "0.679 0.636 0.428 0.507 0.797"
$ws.Range("D59").Value = "0.618 0.525 0.371 0.477 0.765 "
$ws.Range("E59").Value = "0.549 0.485 0.275 0.342 0.667 "
$ws.Range("F59").Value = "0.805 0.859 0.919 0.944 0.977 "
$ws.Range("C60").Value = "0.516 0.201 0.052 0.100 0.684"
$ws.Range("D60").Value = "0.475 0.561 0.238 0.258 0.847 "
$ws.Range("E60").Value = "0.367 0.112 0.026 0.053 0.521 "
$ws.Range("F60").Value = "0.744 0.856 0.928 0.940 0.976 "
$ws.Range("C61").Value = "0.649 0.333 0.200 0.282 0.768"
$ws.Range("D61").Value = "0.676 0.854 0.955 0.833 0.948 "
$ws.Range("E61").Value = "0.501 0.200 0.111 0.164 0.624"
$ws.Range("F61").Value = "0.817 0.877 0.939 0.953 0.983"

# Make the Ubuntu sheet the active tab, matching the selection left after entry
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D61").Select()
